## Journée du 5 mai 2023 — add final "CONCLUSION" slide (slide 16)
## Evolutions des programmes .ino - Présentation pour la journée du 5 mai

$p = $ppt.ActivePresentation

# Append a new slide at the end, using the "Title and Content" layout
# (slideLayout2.xml / ppLayout index 2), same layout used by the slide
# just before it (slide 15 - "Les images avec ESP32 CAM").
$newIndex = $p.Slides.Count + 1
$s = $p.Slides.Add($newIndex, 2)

# --- Title placeholder -------------------------------------------------
$title = $s.Shapes.Item(1).TextFrame.TextRange
$title.Text = "CONCLUSION"
$title.Font.Color.RGB = 12611584   # 0x0070C0 (R=0x00,G=0x70,B=0xC0 -> R+G*256+B*65536)

# --- Content placeholder -------------------------------------------------
$body = $s.Shapes.Item(2).TextFrame.TextRange
$body.Text = "Partage des avis des participants sur l’évolution de l’Intelligence Artificielle dans le futur.`r`rQuelle place aura l’intelligence artificielle dans les prochaines années :`rA l’école`rA la maison`rAu travail`r`rMerci de votre présence !!!`r"

# Paragraphs 4 to 6 ("A l'école", "A la maison", "Au travail") and the
# blank line after them sit at outline level 3 (lvl="2").
$body.Paragraphs(4,1).IndentLevel = 3
$body.Paragraphs(5,1).IndentLevel = 3
$body.Paragraphs(6,1).IndentLevel = 3
$body.Paragraphs(7,1).IndentLevel = 3

# "Merci de votre présence !!!" is large and red.
$merci = $body.Paragraphs(8,1)
$merci.Font.Size = 32
$merci.Font.Color.RGB = 192   # 0xC00000 (R=0xC0,G=0x00,B=0x00 -> R+G*256+B*65536)
